$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Maturity_At_Age: insert a "Sex" column (introduces shared-string "Sex"
#    FIRST so that it lands at shared-string index 59, matching the target).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Maturity_At_Age")
$ws3.Columns("C:C").Insert()
$ws3.Range("C1").Value = "Sex"
$ws3.Range("C2").Value = 1

# Duplicate row 2 into row 3 (same maturity-at-age schedule for both sexes),
# then mark the new row as Sex = 2.
$ws3.Range("A3:AG3").Value = $ws3.Range("A2:AG2").Value()
$ws3.Range("C3").Value = 2

$ws3.Range("A3").Select()
$ws3.Range("A3:XFD3").Select()

# ---------------------------------------------------------------------------
# 2) Weight_At_Age: same column insertion, but the new (male) row has its own
#    distinct weight-at-age values.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Weight_At_Age")
$ws4.Columns("C:C").Insert()
$ws4.Range("C1").Value = "Sex"
$ws4.Range("C2").Value = 1

$ws4.Range("A3").Value = 1
$ws4.Range("B3").Value = "Time_Inv"
$ws4.Range("C3").Value = 2

$maleWeights = @(1.1085,1.4285000000000001,1.7228000000000001,1.9837,2.2088999999999999,2.3995000000000002,2.5586000000000002,2.6899000000000002,2.7974000000000001,2.8847999999999998,2.9554999999999998,3.0125000000000002,3.0583999999999998,3.0951,3.1244999999999998,3.1480000000000001,3.1667999999999998,3.1817000000000002,3.1936,3.2031000000000001,3.2107000000000001,3.2166999999999999,3.2214999999999998,3.2252999999999998,3.2282999999999999,3.2307000000000001,3.2326000000000001,3.2341000000000002,3.2353000000000001,3.2381000000000002)
for ($i = 0; $i -lt $maleWeights.Length; $i++) {
    $ws4.Cells.Item(3, 4 + $i).Value = $maleWeights[$i]
}

$ws4.Range("A3").Select()
$ws4.Range("A3:XFD3").Select()

# ---------------------------------------------------------------------------
# 3) Controls: update simulation/year/sex counts, rename n_fleets ->
#    n_fish_fleets (reusing its row), and append a new survey-fleets row.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Controls")
$ws1.Range("B2").Value = 5
$ws1.Range("B3").Value = 200
$ws1.Range("B5").Value = 2

$ws1.Range("A6").Value = "n_fish_fleets"
$ws1.Range("B6").Value = 2

$ws1.Range("C7").Value = "Number of survey fleets"
$ws1.Range("A7").Value = "n_srv_fleets"
$ws1.Range("B7").Value = 2

$ws1.Range("B5").Select()

# ---------------------------------------------------------------------------
# 4) Age_Bins: selection only (no data change).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Age_Bins")
$ws2.Range("A5:A31").Select()

# ---------------------------------------------------------------------------
# 5) Recruitment_Mortality: bump sigma_rec, select B6, and activate this
#    sheet last so it becomes the workbook's active/visible tab.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Recruitment_Mortality")
$ws5.Range("B5").Value = 1.2

$ws5.Activate()
$ws5.Range("B6").Select()
